$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1, J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style from H1 (existing bold/bordered header style) to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for I2:J66 (rows 2-66), columns I (I0) and J (IF)
$arr = New-Object 'object[,]' 65,2
$arr[0,0] = 9; $arr[0,1] = 9
$arr[1,0] = 9; $arr[1,1] = 9
$arr[2,0] = 8; $arr[2,1] = 8
$arr[3,0] = 7; $arr[3,1] = 8
$arr[4,0] = 7; $arr[4,1] = 8
$arr[5,0] = 7; $arr[5,1] = 7
$arr[6,0] = 9; $arr[6,1] = 10
$arr[7,0] = 8; $arr[7,1] = 8
$arr[8,0] = 8; $arr[8,1] = 9
$arr[9,0] = 8; $arr[9,1] = 8
$arr[10,0] = 8; $arr[10,1] = 9
$arr[11,0] = 9; $arr[11,1] = 10
$arr[12,0] = 9; $arr[12,1] = 9
$arr[13,0] = 8; $arr[13,1] = 8
$arr[14,0] = 6; $arr[14,1] = 6
$arr[15,0] = 8; $arr[15,1] = 9
$arr[16,0] = 6; $arr[16,1] = 6
$arr[17,0] = 6; $arr[17,1] = 6
$arr[18,0] = 8; $arr[18,1] = 8
$arr[19,0] = 7; $arr[19,1] = 7
$arr[20,0] = 7; $arr[20,1] = 7
$arr[21,0] = 6; $arr[21,1] = 7
$arr[22,0] = 9; $arr[22,1] = 9
$arr[23,0] = 9; $arr[23,1] = 9
$arr[24,0] = 8; $arr[24,1] = 8
$arr[25,0] = 9; $arr[25,1] = 9
$arr[26,0] = 7; $arr[26,1] = 7
$arr[27,0] = 6; $arr[27,1] = 7
$arr[28,0] = 7; $arr[28,1] = 8
$arr[29,0] = 8; $arr[29,1] = 8
$arr[30,0] = 6; $arr[30,1] = 6
$arr[31,0] = 5; $arr[31,1] = 6
$arr[32,0] = 8; $arr[32,1] = 9
$arr[33,0] = 8; $arr[33,1] = 9
$arr[34,0] = 4; $arr[34,1] = 6
$arr[35,0] = 10; $arr[35,1] = 10
$arr[36,0] = 8; $arr[36,1] = 8
$arr[37,0] = 9; $arr[37,1] = 9
$arr[38,0] = 5; $arr[38,1] = 6
$arr[39,0] = 4; $arr[39,1] = 4
$arr[40,0] = 7; $arr[40,1] = 8
$arr[41,0] = 7; $arr[41,1] = 7
$arr[42,0] = 9; $arr[42,1] = 9
$arr[43,0] = 6; $arr[43,1] = 6
$arr[44,0] = 6; $arr[44,1] = 6
$arr[45,0] = 7; $arr[45,1] = 7
$arr[46,0] = 5; $arr[46,1] = 6
$arr[47,0] = 7; $arr[47,1] = 7
$arr[48,0] = 8; $arr[48,1] = 9
$arr[49,0] = 6; $arr[49,1] = 7
$arr[50,0] = 7; $arr[50,1] = 7
$arr[51,0] = 8; $arr[51,1] = 9
$arr[52,0] = 6; $arr[52,1] = 7
$arr[53,0] = 7; $arr[53,1] = 7
$arr[54,0] = 6; $arr[54,1] = 7
$arr[55,0] = 7; $arr[55,1] = 8
$arr[56,0] = 7; $arr[56,1] = 8
$arr[57,0] = 5; $arr[57,1] = 6
$arr[58,0] = 7; $arr[58,1] = 7
$arr[59,0] = 7; $arr[59,1] = 8
$arr[60,0] = 7; $arr[60,1] = 7
$arr[61,0] = 9; $arr[61,1] = 9
$arr[62,0] = 6; $arr[62,1] = 6
$arr[63,0] = 7; $arr[63,1] = 7
$arr[64,0] = 7; $arr[64,1] = 7

$ws.Range("I2:J66").Value = $arr